# Generate Report for Handoff
# Adds a new tracked file (4c04c876-246e-47a4-8d39-0350ea6ec30b.md) to the
# localization-status workbook, inserted just above the always-last
# ".localization-config" row on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (File Name / zh-cn / de-de) -------------------------
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Rows.Item(7).Insert()

$ov.Range("A7").Value = "4c04c876-246e-47a4-8d39-0350ea6ec30b.md"
$ov.Range("B7").Value = "Ready for handoff"
$ov.Range("C7").Value = "Ready for handoff"

# Row 8 keeps the ".localization-config" content that used to live in row 7
# (Excel's row-insert already shifted the cell text down); just make sure
# it is the expected value.
$ov.Range("A8").Value = ".localization-config"
$ov.Range("B8").Value = "Not to be localized"
$ov.Range("C8").Value = "Not to be localized"

# Hyperlink refs don't auto-shift with a row insert in this engine, so
# rebuild the hyperlinks collection for the sheet from scratch.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32a840d272592c4370bfbaaa97c8e4441f143c3d/e2e/376d69e4-644e-4015-b20f-f4155f460ced.md", "", "", "376d69e4-644e-4015-b20f-f4155f460ced.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/e2e/4b874caa-3a1c-443d-8c7a-18d54918465a.md", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/87764c028ea30d400f5be49b855fe65bffb2936d/e2e/89b04cff-4acc-4251-bf2b-f0c09418649e.md", "", "", "89b04cff-4acc-4251-bf2b-f0c09418649e.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/87764c028ea30d400f5be49b855fe65bffb2936d/e2e/93a5442b-4b8c-4207-a35a-7ed7a4b2544a.md", "", "", "93a5442b-4b8c-4207-a35a-7ed7a4b2544a.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/179f4c8669219b6240d9c018b1cd6fe0647fffeb/e2e/244e4e48-50cf-4170-8608-4cd025acf3f9.md", "", "", "244e4e48-50cf-4170-8608-4cd025acf3f9.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/e2e/4c04c876-246e-47a4-8d39-0350ea6ec30b.md", "", "", "4c04c876-246e-47a4-8d39-0350ea6ec30b.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" ----------------------------------------------------------
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Rows.Item(7).Insert()

$zh.Range("A7").Value = "4c04c876-246e-47a4-8d39-0350ea6ec30b.md"
$zh.Range("B7").Value = "Ready for handoff"
$zh.Range("C7").Value = "4c04c876-246e-47a4-8d39-0350ea6ec30b.276a943e2f764fdfb57aa1437f1267e6e923b273.zh-cn.xlf"
$zh.Range("D7").Value = "2016-02-22 04:19:26"
$zh.Range("G7").Value = "0001-01-01 00:00:00"
$zh.Range("H7").Value = "Include"

$zh.Range("A8").Value = ".localization-config"
$zh.Range("B8").Value = "Not to be localized"
$zh.Range("D8").Value = "0001-01-01 00:00:00"
$zh.Range("G8").Value = "0001-01-01 00:00:00"
$zh.Range("H8").Value = "Ignored"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32a840d272592c4370bfbaaa97c8e4441f143c3d/e2e/376d69e4-644e-4015-b20f-f4155f460ced.md", "", "", "376d69e4-644e-4015-b20f-f4155f460ced.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b07e1a6fc504d205291596f8d91ce28d721e370/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/376d69e4-644e-4015-b20f-f4155f460ced.3bfcb2ad13ebb5475a745cba1d66ecb4b9771ced.zh-cn.xlf", "", "", "376d69e4-644e-4015-b20f-f4155f460ced.3bfcb2ad13ebb5475a745cba1d66ecb4b9771ced.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/e2e/4b874caa-3a1c-443d-8c7a-18d54918465a.md", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9f50c63a9c20e4e60c3f3eca78396f289c58db0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.zh-cn.xlf", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b2e4c001a477e9ab0763448d0b275267558e7561/e2e/4b874caa-3a1c-443d-8c7a-18d54918465a.md", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d64d81dbaaefce7a3726c9ab6ccc241066ab84aa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.zh-cn.xlf", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/87764c028ea30d400f5be49b855fe65bffb2936d/e2e/89b04cff-4acc-4251-bf2b-f0c09418649e.md", "", "", "89b04cff-4acc-4251-bf2b-f0c09418649e.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c856d36f2d69ca8a0da540f1baf8ef6aa978944/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/89b04cff-4acc-4251-bf2b-f0c09418649e.83327157723093f3f5d1c8b535f67b357fb0cdc9.zh-cn.xlf", "", "", "89b04cff-4acc-4251-bf2b-f0c09418649e.83327157723093f3f5d1c8b535f67b357fb0cdc9.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/87764c028ea30d400f5be49b855fe65bffb2936d/e2e/93a5442b-4b8c-4207-a35a-7ed7a4b2544a.md", "", "", "93a5442b-4b8c-4207-a35a-7ed7a4b2544a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c856d36f2d69ca8a0da540f1baf8ef6aa978944/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/93a5442b-4b8c-4207-a35a-7ed7a4b2544a.11169a5bf040afb60722664672641e052b89f185.zh-cn.xlf", "", "", "93a5442b-4b8c-4207-a35a-7ed7a4b2544a.11169a5bf040afb60722664672641e052b89f185.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/179f4c8669219b6240d9c018b1cd6fe0647fffeb/e2e/244e4e48-50cf-4170-8608-4cd025acf3f9.md", "", "", "244e4e48-50cf-4170-8608-4cd025acf3f9.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b6937cb0f6eaffb519106b1d6f9b162b16c7693a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/244e4e48-50cf-4170-8608-4cd025acf3f9.7a58c17f6c0d9ed2540461573228c309ea054904.zh-cn.xlf", "", "", "244e4e48-50cf-4170-8608-4cd025acf3f9.7a58c17f6c0d9ed2540461573228c309ea054904.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/e2e/4c04c876-246e-47a4-8d39-0350ea6ec30b.md", "", "", "4c04c876-246e-47a4-8d39-0350ea6ec30b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/276a943e2f764fdfb57aa1437f1267e6e923b273/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4c04c876-246e-47a4-8d39-0350ea6ec30b.276a943e2f764fdfb57aa1437f1267e6e923b273.zh-cn.xlf", "", "", "4c04c876-246e-47a4-8d39-0350ea6ec30b.276a943e2f764fdfb57aa1437f1267e6e923b273.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" ----------------------------------------------------------
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Rows.Item(7).Insert()

$de.Range("A7").Value = "4c04c876-246e-47a4-8d39-0350ea6ec30b.md"
$de.Range("B7").Value = "Ready for handoff"
$de.Range("C7").Value = "4c04c876-246e-47a4-8d39-0350ea6ec30b.276a943e2f764fdfb57aa1437f1267e6e923b273.de-de.xlf"
$de.Range("D7").Value = "2016-02-22 04:19:39"
$de.Range("G7").Value = "0001-01-01 00:00:00"
$de.Range("H7").Value = "Include"

$de.Range("A8").Value = ".localization-config"
$de.Range("B8").Value = "Not to be localized"
$de.Range("D8").Value = "0001-01-01 00:00:00"
$de.Range("G8").Value = "0001-01-01 00:00:00"
$de.Range("H8").Value = "Ignored"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32a840d272592c4370bfbaaa97c8e4441f143c3d/e2e/376d69e4-644e-4015-b20f-f4155f460ced.md", "", "", "376d69e4-644e-4015-b20f-f4155f460ced.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1fac8a1f520450b07d964d4112ffd927b37fb5e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/376d69e4-644e-4015-b20f-f4155f460ced.3bfcb2ad13ebb5475a745cba1d66ecb4b9771ced.de-de.xlf", "", "", "376d69e4-644e-4015-b20f-f4155f460ced.3bfcb2ad13ebb5475a745cba1d66ecb4b9771ced.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/e2e/4b874caa-3a1c-443d-8c7a-18d54918465a.md", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7377ee667c33bd03a0c7586950c6cdf9f142b7d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.de-de.xlf", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/01e83f3c4debe2dd23a9447609c7d84b8af85ca5/e2e/4b874caa-3a1c-443d-8c7a-18d54918465a.md", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/41a78e90b24770d131681d1927139db3f39a7688/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.de-de.xlf", "", "", "4b874caa-3a1c-443d-8c7a-18d54918465a.05d3c026fc73a930ab607897c5660517dd1cd4db.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/87764c028ea30d400f5be49b855fe65bffb2936d/e2e/89b04cff-4acc-4251-bf2b-f0c09418649e.md", "", "", "89b04cff-4acc-4251-bf2b-f0c09418649e.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13921c14f0e896eb7eabf28cf8bf6f44773c08ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/89b04cff-4acc-4251-bf2b-f0c09418649e.83327157723093f3f5d1c8b535f67b357fb0cdc9.de-de.xlf", "", "", "89b04cff-4acc-4251-bf2b-f0c09418649e.83327157723093f3f5d1c8b535f67b357fb0cdc9.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/87764c028ea30d400f5be49b855fe65bffb2936d/e2e/93a5442b-4b8c-4207-a35a-7ed7a4b2544a.md", "", "", "93a5442b-4b8c-4207-a35a-7ed7a4b2544a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13921c14f0e896eb7eabf28cf8bf6f44773c08ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/93a5442b-4b8c-4207-a35a-7ed7a4b2544a.11169a5bf040afb60722664672641e052b89f185.de-de.xlf", "", "", "93a5442b-4b8c-4207-a35a-7ed7a4b2544a.11169a5bf040afb60722664672641e052b89f185.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/179f4c8669219b6240d9c018b1cd6fe0647fffeb/e2e/244e4e48-50cf-4170-8608-4cd025acf3f9.md", "", "", "244e4e48-50cf-4170-8608-4cd025acf3f9.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01ddab9e2911c5246b2ca1651a829874570aff41/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/244e4e48-50cf-4170-8608-4cd025acf3f9.7a58c17f6c0d9ed2540461573228c309ea054904.de-de.xlf", "", "", "244e4e48-50cf-4170-8608-4cd025acf3f9.7a58c17f6c0d9ed2540461573228c309ea054904.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/e2e/4c04c876-246e-47a4-8d39-0350ea6ec30b.md", "", "", "4c04c876-246e-47a4-8d39-0350ea6ec30b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/276a943e2f764fdfb57aa1437f1267e6e923b273/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4c04c876-246e-47a4-8d39-0350ea6ec30b.276a943e2f764fdfb57aa1437f1267e6e923b273.de-de.xlf", "", "", "4c04c876-246e-47a4-8d39-0350ea6ec30b.276a943e2f764fdfb57aa1437f1267e6e923b273.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/b202bf120df79d0f720fc89561e60fc17e818cf1/.localization-config", "", "", ".localization-config") | Out-Null

"Generate Report for Handoff: applied"
